# "notes from july 4" - add a new section to the end of the presentation:
#   - a bold heading "Identification of models with uncertainty, learning,
#     and human capital sorting"
#   - four new bulleted (ListParagraph) points under it, reusing the
#     existing trailing empty bullet paragraph for the last point.

$d = $word.ActiveDocument

function New-DocXmlPackage([string]$bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1) New bold heading paragraph (plain paragraph, no list style).
$headingP = '<w:p>' +
              '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr>' +
                '<w:t>Identification of models with uncertainty, learning, and human capital sorting</w:t>' +
              '</w:r>' +
            '</w:p>'
$target = $d.Paragraphs.Last.Range
$null = $target.InsertXML((New-DocXmlPackage $headingP))

# 2) "Tldr: ..." bullet, with the spell-check markers around "Tldr".
$tldrP = '<w:p>' +
           '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:t>Tldr</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve">: if expected to learn a lot of human capital, can be paid a lower wage </w:t></w:r>' +
         '</w:p>'
$target = $d.Paragraphs.Last.Range
$null = $target.InsertXML((New-DocXmlPackage $tldrP))

# 3) "Sorting matters for wage inequality" bullet.
$sortingP = '<w:p>' +
              '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
              '<w:r><w:t>Sorting matters for wage inequality</w:t></w:r>' +
            '</w:p>'
$target = $d.Paragraphs.Last.Range
$null = $target.InsertXML((New-DocXmlPackage $sortingP))

# 4) "Compare between 1st and 2nd best choices ..." bullet, with superscript ordinals.
$compareP = '<w:p>' +
              '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
              '<w:r><w:t>Compare between 1</w:t></w:r>' +
              '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> and 2</w:t></w:r>' +
              '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> best choices and the difference vs. compensating differential</w:t></w:r>' +
            '</w:p>'
$target = $d.Paragraphs.Last.Range
$null = $target.InsertXML((New-DocXmlPackage $compareP))

# 5) "AKM underestimate impact of sorting" - goes directly into the original
#    trailing empty bullet paragraph (so no extra empty paragraph is left
#    behind at the end of the document).
$target = $d.Paragraphs.Last.Range
$null = $target.InsertBefore("AKM underestimate impact of sorting")
